# Reducing line spacing in title
#
# 1) Bump the Title style's paragraph line spacing to 1.5 lines
#    (renders as <w:spacing w:after="300" w:line="360" w:lineRule="auto"/>).
# 2) Word re-anchors the hidden "_GoBack" bookmark (last-edit marker) to
#    wherever the edit actually happened - here, the very start of the
#    Title paragraph, which ends up splitting the "Title" run into "T" /
#    "itle" around the bookmark.

$d = $word.ActiveDocument

$titleStyle = $d.Styles("Title")
$titleStyle.ParagraphFormat.LineSpacingRule = 1   # wdLineSpace1pt5

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$editSpot = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $editSpot)
